# CreateProject.xlsx - 30-11 updated for Created Project TAsked
#
# Summary of changes:
#  - NewProject: remove the two extra sample rows (rows 3-5), and update the
#    remaining record's ProjectCode/ProjectName ("YT626" -> "YTP565") and
#    ProjectManager ("Alok Malviya" -> "Demo Titan").
#  - NewTask: update the task record (Task06 -> Task08, StartDate/DueDate
#    values swapped) and make this sheet the active one.
#  - ProjectTeam: remove the TeamMember/Role columns, keeping only
#    EngageDate/ExRelasedate, and move the selection off this sheet.

$wb = $excel.ActiveWorkbook

$wsProject = $wb.Worksheets.Item("NewProject")
$wsTask    = $wb.Worksheets.Item("NewTask")
$wsTeam    = $wb.Worksheets.Item("ProjectTeam")

# --- NewProject sheet -------------------------------------------------
# Drop the extra sample rows (3,4) and the stray formatted row (5),
# keeping just the header + one data row.
$wsProject.Rows("3:5").Delete()

# Update the remaining data row with the new project code/name and manager.
$wsProject.Range("A2").Value = "YTP565"
$wsProject.Range("B2").Value = "YTP565"
$wsProject.Range("K2").Value = "Demo Titan"

$wsProject.Range("A1:XFD1048576").Select()

# --- NewTask sheet ------------------------------------------------------
$wsTask.Range("A2").Value = "Task08"
$wsTask.Range("C2").Value = "1"
$wsTask.Range("D2").Value = "3"

$wsTask.Range("D9").Select()
$wsTask.Activate()

# --- ProjectTeam sheet ---------------------------------------------------
# Remove the TeamMember and Role columns entirely; EngageDate/ExRelasedate
# shift left into columns A and B.
$wsTeam.Columns("A:B").Delete()

$wsTeam.Range("G7").Select()
